# Update Name of Algo
# Apply updated KNN imputation results to specific cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.475999999999999
$ws.Range("B9").Value = 6.485000000000001
$ws.Range("D11").Value = -8.316999999999998
$ws.Range("B18").Value = 6.351
$ws.Range("B20").Value = 6.667999999999999
